# Update "want to go" counts (column F) for two events that each appear
# on the "展览" sheet and the "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 285
$wsExhibit.Range("F9").Value = 2017

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F6").Value = 285
$wsAll.Range("F13").Value = 2017
